$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(2901, 3591, 4028, 4338, 4586, 4794, 4794, 4839, 4852, 4852, 4852, 4852, 4852, 4897)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
